# Add a new meeting-diary entry for the Oct 27 discussion (row 25),
# matching the formatting of the preceding entries (rows 21-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the style of the row above (row 24) first, so the new row matches the
# formatting of the rest of the table (date/time number formats, wrap text,
# quote-prefix on the Discussions cell).
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Date / start-time / end-time for the new meeting
$ws.Range("A25").Value = (Get-Date -Year 2023 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B25").Value = 0.47916666666666669
$ws.Range("C25").Value = 0.625

# Members present (same group as the previous few entries)
$ws.Range("D25").Value = "Aishwarya Anil Kumar (32644329) / Chih Hui Wang (33209006) / Shreyansh Mahtolia (33509115)"

# Discussion notes for the new meeting (leading apostrophe forces Excel's
# "text" quote-prefix interpretation, matching the other Discussions cells
# whose text also starts with "-")
$ws.Range("E25").Value = "'- Discuss and modify the final report version.
- Make minor corrections.
- Discuss video presentation submission."

# Match row height of the adjacent rows (wrap text etc.)
$ws.Rows("25").RowHeight = 51

# Update the active selection to C26, as in the saved workbook
$ws.Range("C26").Select()
